# Adjusted eval parameters, added loader for external evaluation.
# Rewrites the A1:P16 evaluation-parameter grid on the active worksheet
# with the updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(37, 1, 42, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 13),
    @(8, 46, 15, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 4),
    @(11, 1, 67, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 74),
    @(0, 0, 0, 72, 0, 8, 0, 0, 0, 0, 0, 0, 0, 0, 0, 15),
    @(0, 0, 0, 0, 61, 13, 0, 0, 0, 0, 0, 0, 0, 0, 0, 9),
    @(0, 0, 0, 0, 0, 77, 0, 0, 0, 0, 0, 0, 0, 0, 0, 160),
    @(0, 0, 0, 0, 0, 0, 68, 1, 11, 0, 0, 0, 0, 0, 0, 16),
    @(0, 0, 0, 0, 0, 0, 1, 63, 10, 0, 0, 0, 0, 0, 0, 14),
    @(0, 0, 0, 0, 0, 0, 2, 1, 77, 0, 0, 0, 0, 0, 0, 127),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 72, 0, 8, 0, 0, 0, 8),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 69, 6, 0, 0, 0, 10),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 9, 0, 69, 0, 0, 0, 114),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 70, 0, 10, 11),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 62, 5, 12),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 78, 172),
    @(0, 11, 1, 0, 6, 3, 0, 6, 0, 0, 4, 2, 1, 13, 2, 0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}
